{"js": "// The document contains two adjacent \"JOB\"/\"Job\" heading blocks (an\n// Arabic-language resource). The first block is:\n//   - Heading2 paragraph: \"JOB\"\n//   - Normal paragraph containing only an italic run reading \"Job\"\n//   - Normal paragraph containing a single space\n// The edit removes the whole italic-only \"Job\" paragraph that sits\n// between the \"JOB\" heading and the paragraph with the lone space.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load the text for every paragraph so we can unambiguously find the\n// specific empty/italic \"Job\" paragraph that immediately follows the\n// \"JOB\" Heading 2 (there is a second, unrelated \"Job\" heading further\n// down in the document with real body text beneath it - that one must\n// stay untouched).\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text,style\");\n}\nawait context.sync();\n\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"Heading 2\" && para.text.trim() === \"JOB\") {\n    const next = paragraphs.items[i + 1];\n    if (next.style === \"Normal\" && next.text.trim() === \"Job\") {\n      targetIndex = i + 1;\n      break;\n    }\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate the italic 'Job' paragraph to delete.\");\n}\n\nparagraphs.items[targetIndex].delete();\nawait context.sync();\n", "ps1": "# The document contains two adjacent \"JOB\"/\"Job\" heading blocks (an\n# Arabic-language resource). The first block is:\n#   - Heading 2 paragraph: \"JOB\"\n#   - Normal paragraph containing only an italic run reading \"Job\"\n#   - Normal paragraph containing a single space\n# The edit removes the whole italic-only \"Job\" paragraph that sits\n# between the \"JOB\" heading and the paragraph with the lone space.\n#\n# We locate it by scanning paragraphs for a Heading 2 \"JOB\" immediately\n# followed by a Normal paragraph whose text is \"Job\" - this disambiguates\n# it from the unrelated, second \"Job\" Heading 2 further down the document\n# (which has real body text under it and must stay untouched).\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n$targetIndex = -1\nfor ($i = 1; $i -lt $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Range.Style.NameLocal\n    $txt = $p.Range.Text.Trim()\n    if ($styleName -eq \"Heading 2\" -and $txt -eq \"JOB\") {\n        $next = $d.Paragraphs.Item($i + 1)\n        $nextStyle = $next.Range.Style.NameLocal\n        $nextTxt = $next.Range.Text.Trim()\n        if ($nextStyle -eq \"Normal\" -and $nextTxt -eq \"Job\") {\n            $targetIndex = $i + 1\n            break\n        }\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the italic 'Job' paragraph to delete.\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.Delete()\n"}
